$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (no date-like auto-conversion risk) ---
$ws.Range("B1").Value = "akhil"
$ws.Range("B2").Value = "rk indusreis"
$ws.Range("B7").Value = "31-10-2020"
$ws.Range("B10").Value = "gst report test"

# --- Numeric cell ---
$ws.Range("E10").Value = 48.14

# --- Date-like text cells that Excel would otherwise auto-convert to date serials ---
# Use a scratch area (column H, far outside the used range) to stage the
# text-formatted value, then copy only the value into the target cell so the
# target cell's existing style/border (if any) is left completely untouched.

# B5: 2020-07-29 -> 2020-10-05
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "2020-10-05"
$ws.Range("H1").Copy()
$ws.Range("B5").PasteSpecial(-4163)  # xlPasteValues

# B6: 2020-07-16 -> 01-10-2020
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "01-10-2020"
$ws.Range("H2").Copy()
$ws.Range("B6").PasteSpecial(-4163)  # xlPasteValues

# A10: 2020-07-29 -> 2020-10-05 (preserve its existing style/border via stash+restore)
$ws.Range("A10").Copy()
$ws.Range("H4").PasteSpecial(-4122)  # xlPasteFormats (stash original format)

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2020-10-05"
$ws.Range("H3").Copy()
$ws.Range("A10").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("H4").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats (restore original format)

# clean up scratch cells
$ws.Range("H1:H4").Clear()
